$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '37.436.20'
$ws.Range('E2').Value = '  +1.40%  '
$ws.Range('D3').Value = '2.026.31'
$ws.Range('E3').Value = '  +2.08%  '
$ws.Range('E4').Value = '  -0.02%  '
$ws.Range('D5').Value = '255.16'
$ws.Range('E5').Value = '  +4.08%  '
$ws.Range('D6').Value = '0.621'
$ws.Range('E6').Value = '  -1.22%  '
$ws.Range('D7').Value = '1.00'
$ws.Range('E7').Value = '  -0.03%  '
$ws.Range('D8').Value = '57.20'
$ws.Range('E8').Value = '  -5.93%  '
$ws.Range('D9').Value = '0.385'
$ws.Range('E9').Value = '  +0.61%  '
$ws.Range('D10').Value = '0.0789'
$ws.Range('E10').Value = '  -1.20%  '
$ws.Range('E11').Value = '  -1.75%  '
$ws.Range('D12').Value = '14.71'
$ws.Range('E12').Value = '  -2.28%  '
$ws.Range('D13').Value = '2.327.41'
$ws.Range('E13').Value = '  +2.24%  '
$ws.Range('D14').Value = '0.816'
$ws.Range('E14').Value = '  -3.29%  '
$ws.Range('D15').Value = '21.23'
$ws.Range('E15').Value = '  -3.84%  '
$ws.Range('D16').Value = '5.33'
$ws.Range('E16').Value = '  -2.37%  '
$ws.Range('D17').Value = '2.031.50'
$ws.Range('E17').Value = '  +2.25%  '
$ws.Range('D18').Value = '37.406.22'
$ws.Range('E18').Value = '  +1.54%  '
$ws.Range('D19').Value = '69.68'
$ws.Range('E19').Value = '  -0.75%  '
$ws.Range('D20').Value = '0.0₃0851'
$ws.Range('E20').Value = '  -1.17%  '
$ws.Range('D21').Value = '5.18'
$ws.Range('E21').Value = '  +0.04%  '
$ws.Range('D22').Value = '228.40'
$ws.Range('E22').Value = '  -0.70%  '
$ws.Range('D23').Value = '2.62'
$ws.Range('E23').Value = '  +3.98%  '
$ws.Range('E24').Value = '  +0.02%  '
$ws.Range('E25').Value = '  -0.97%  '
$ws.Range('B26').Value = 'Monero'
$ws.Range('C26').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D26').Value = '162.62'
$ws.Range('E26').Value = '  -0.26%  '
$ws.Range('B27').Value = 'Cosmos'
$ws.Range('C27').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D27').Value = '9.03'
$ws.Range('E27').Value = '  -2.58%  '
$ws.Range('B28').Value = 'EthereumClassic'
$ws.Range('C28').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D28').Value = '19.86'
$ws.Range('E28').Value = '  +1.56%  '
$ws.Range('B29').Value = 'Kaspa'
$ws.Range('C29').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D29').Value = '0.131'
$ws.Range('E29').Value = '  -12.05%  '
$ws.Range('D30').Value = '1.35'
$ws.Range('E30').Value = '  -1.89%  '
$ws.Range('E31').Value = '  -1.43%  '
$ws.Range('D32').Value = '0.0662'
$ws.Range('E32').Value = '  +6.66%  '
$ws.Range('D33').Value = '4.70'
$ws.Range('E33').Value = '  -3.56%  '
$ws.Range('D34').Value = '4.52'
$ws.Range('E34').Value = '  -0.47%  '
$ws.Range('E35').Value = '  +6.73%  '
$ws.Range('E36').Value = '  -0.02%  '
$ws.Range('E37').Value = '  +1.93%  '
$ws.Range('D38').Value = '3.37'
$ws.Range('E38').Value = '  +0.93%  '
$ws.Range('D39').Value = '5.34'
$ws.Range('E39').Value = '  -3.49%  '
$ws.Range('E40').Value = '  +3.81%  '
$ws.Range('D41').Value = '0.0965'
$ws.Range('E41').Value = '  -3.14%  '
$ws.Range('B42').Value = 'VeChain'
$ws.Range('C42').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D42').Value = '0.0215'
$ws.Range('E42').Value = '  +1.06%  '
$ws.Range('B43').Value = 'TrustWalletToken'
$ws.Range('C43').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D43').Value = '1.18'
$ws.Range('E43').Value = '  +0.25%  '
$ws.Range('D44').Value = '1.394.89'
$ws.Range('E44').Value = '  +1.85%  '
$ws.Range('D45').Value = '16.00'
$ws.Range('E45').Value = '  -2.93%  '
$ws.Range('D46').Value = '90.42'
$ws.Range('E46').Value = '  +0.28%  '
$ws.Range('D47').Value = '1.04'
$ws.Range('E47').Value = '  +0.14%  '
$ws.Range('D48').Value = '7.33'
$ws.Range('E48').Value = '  +1.00%  '
$ws.Range('E49').Value = '  +1.87%  '
$ws.Range('D50').Value = '2.01'
$ws.Range('E50').Value = '  +1.20%  '
$ws.Range('D51').Value = '2.219.38'
$ws.Range('E51').Value = '  +2.32%  '
